# The document has two image logos (a Pearson PNG logo in the footers and a
# BTEC JPG logo in the headers), each placed twice (once in the "primary"
# header/footer and once in the "first page" header/footer of the only
# section). The edit swaps the part names used for these inline pictures:
#   - the Pearson PNG logos go from "image2.png" -> "image1.png"
#   - the BTEC JPG logos go from "image1.jpg" -> "image2.jpg"
#
# InlineShape has no writable Name property in the Word object model, but
# Shape does, and <InlineShape>.ConvertToShape() / <Shape>.ConvertToInlineShape()
# let us round-trip through a floating shape just long enough to rename it
# while keeping it as an inline picture afterwards.

function Rename-LogoInRange($range, $newName) {
    $inline = $range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

$d = $word.ActiveDocument
$sec = $d.Sections.First

# Footers: Pearson logo, image2.png -> image1.png
Rename-LogoInRange $sec.Footers.Item(1).Range "image1.png"
Rename-LogoInRange $sec.Footers.Item(2).Range "image1.png"

# Headers: BTEC logo, image1.jpg -> image2.jpg
Rename-LogoInRange $sec.Headers.Item(1).Range "image2.jpg"
Rename-LogoInRange $sec.Headers.Item(2).Range "image2.jpg"
